$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("E1").Value = "numerator"
$ws.Range("F1").Value = "numerator_desc"
$ws.Range("G1").Value = "follow_up"

# Update data rows 2-5
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 5).Value = $false   # E: numerator (boolean)
    $ws.Cells.Item($r, 6).Value = "No screening recorded"  # F: numerator_desc
    $ws.Cells.Item($r, 7).Value = $true    # G: follow_up (boolean)
}

# H column (medicaid) swap for rows 4 and 5
$ws.Cells.Item(4, 8).Value = $true
$ws.Cells.Item(5, 8).Value = $false
